# "Add files via upload" — refresh the player list on Sheet1 with the
# latest "Out of PO" roster. The list grew by one row (Amen Thompson was
# inserted) and several players were re-shuffled to new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Luguentz Dort", "SG,SF", "Oklahoma City Thunder"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Ochai Agbaji", "SG,SF", "Toronto Raptors"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
}
